$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1050.4117
$ws.Range("I38").Value = 79.875
$ws.Range("J38").Value = 1913.1111
$ws.Range("K38").Value = 239.625
$ws.Range("L38").Value = 5739.3333
$ws.Range("M38").Value = 132.375
$ws.Range("N38").Value = -6483.3333

$ws.Range("H129").Value = 1514.3077
$ws.Range("I129").Value = 295.66666
$ws.Range("J129").Value = 1673.2609
$ws.Range("K129").Value = 886.9999799999999
$ws.Range("L129").Value = 5019.7827
$ws.Range("M129").Value = 4113.00002
$ws.Range("N129").Value = -15019.7827

$ws.Range("H137").Value = 1610.4482
$ws.Range("I137").Value = 1388.625
$ws.Range("J137").Value = 1883.4615
$ws.Range("K137").Value = 4165.875
$ws.Range("L137").Value = 5650.3845
$ws.Range("M137").Value = -1615.875
$ws.Range("N137").Value = -10750.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1143.4445
$ws.Range("I4").Value = 1558.2
$ws.Range("J4").Value = 625
$ws.Range("K4").Value = 1558.2
$ws.Range("L4").Value = 625
$ws.Range("M4").Value = -1442.2
$ws.Range("N4").Value = -857

$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H23").Value = 18000
$ws.Range("J23").Value = 18000
$ws.Range("L23").Value = 18000
$ws.Range("N23").Value = -18518

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("M37").ClearContents()

$ws.Range("H44").Value = 9959
$ws.Range("J44").Value = 9959
$ws.Range("L44").Value = 9959
$ws.Range("N44").Value = -10935

$ws.Range("H63").Value = 2950
$ws.Range("I63").Value = 2950
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2950
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = -2264
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 2950
$ws.Range("I66").Value = 2950
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14750
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = -11318
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 17309802
$ws.Range("I74").Value = 23685322
$ws.Range("J74").Value = 4818.857
$ws.Range("K74").Value = 23685322
$ws.Range("L74").Value = 4818.857
$ws.Range("M74").Value = -23684448
$ws.Range("N74").Value = -6566.857

$ws.Range("H77").Value = 17309802
$ws.Range("I77").Value = 23685322
$ws.Range("J77").Value = 4818.857
$ws.Range("K77").Value = 118426610
$ws.Range("L77").Value = 24094.285
$ws.Range("M77").Value = -118422242
$ws.Range("N77").Value = -32830.285

$ws.Range("H80").Value = 27866.666
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 39800
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 39800
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -41796

$ws.Range("H83").Value = 27866.666
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 39800
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 119400
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -129384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 107.333336
$ws.Range("I22").Value = 88.8
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 88.8
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 84.2
$ws.Range("N22").Value = -546

$ws.Range("H82").Value = 20969
$ws.Range("I82").Value = 14950.4
$ws.Range("J82").Value = 31000
$ws.Range("K82").Value = 14950.4
$ws.Range("L82").Value = 31000
$ws.Range("M82").Value = -14567.4
$ws.Range("N82").Value = -31766

$ws.Range("H85").Value = 20969
$ws.Range("I85").Value = 14950.4
$ws.Range("J85").Value = 31000
$ws.Range("K85").Value = 14950.4
$ws.Range("L85").Value = 31000
$ws.Range("M85").Value = -13624.4
$ws.Range("N85").Value = -33652

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("L95").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 6249.75
$ws.Range("I41").Value = 6249.75
$ws.Range("K41").Value = 6249.75
$ws.Range("M41").Value = -5821.75

$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 15000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 15000
$ws.Range("N50").Value = 0
$ws.Range("M50").Value = -14375
$ws.Range("L50").ClearContents()

$ws.Range("H59").Value = 22000
$ws.Range("J59").Value = 22000
$ws.Range("L59").Value = 22000
$ws.Range("N59").Value = -24290

$ws.Range("H60").Value = 12531.107
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 12531.107
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = 12531.107
$ws.Range("N60").Value = -13553.107
$ws.Range("L60").ClearContents()

$ws.Range("H62").Value = 100005150
$ws.Range("I62").Value = 3233.3333
$ws.Range("J62").Value = 142863100
$ws.Range("K62").Value = 3233.3333
$ws.Range("L62").Value = 142863100
$ws.Range("M62").Value = -2609.3333
$ws.Range("N62").Value = -142864348

$ws.Range("H65").Value = 100005150
$ws.Range("I65").Value = 3233.3333
$ws.Range("J65").Value = 142863100
$ws.Range("K65").Value = 16166.6665
$ws.Range("L65").Value = 714315500
$ws.Range("M65").Value = -13046.6665
$ws.Range("N65").Value = -714321740

$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("K74").Value = 25000
$ws.Range("M74").Value = -24126

$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70632

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 266.2
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 395.5
$ws.Range("K17").Value = 540
$ws.Range("L17").Value = 1186.5
$ws.Range("M17").Value = -371
$ws.Range("N17").Value = -1524.5

$ws.Range("H132").Value = 1228.9286
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 1242.0834
$ws.Range("K132").Value = 10350
$ws.Range("L132").Value = 11178.7506
$ws.Range("M132").Value = -7820
$ws.Range("N132").Value = -16238.7506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57.7
$ws.Range("I2").Value = 41.7
$ws.Range("J2").Value = 73.7
$ws.Range("K2").Value = 41.7
$ws.Range("L2").Value = 73.7
$ws.Range("M2").Value = 71.3
$ws.Range("N2").Value = -299.7

$ws.Range("H126").Value = 2087.8484
$ws.Range("I126").Value = 2095.96
$ws.Range("J126").Value = 2062.5
$ws.Range("K126").Value = 6287.88
$ws.Range("L126").Value = 6187.5
$ws.Range("M126").Value = -3817.88
$ws.Range("N126").Value = -11127.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1436.4546
$ws.Range("I46").Value = 966.5
$ws.Range("J46").Value = 2000.4
$ws.Range("K46").Value = 966.5
$ws.Range("L46").Value = 2000.4
$ws.Range("M46").Value = -778.5
$ws.Range("N46").Value = -2376.4

$ws.Range("H68").Value = 14712265
$ws.Range("I68").Value = 84584210
$ws.Range("J68").Value = 2381.7368
$ws.Range("K68").Value = 84584210
$ws.Range("L68").Value = 2381.7368
$ws.Range("M68").Value = -84583461
$ws.Range("N68").Value = -3879.7368

$ws.Range("H71").Value = 14712265
$ws.Range("I71").Value = 84584210
$ws.Range("J71").Value = 2381.7368
$ws.Range("K71").Value = 422921050
$ws.Range("L71").Value = 11908.684
$ws.Range("M71").Value = -422917306
$ws.Range("N71").Value = -19396.684

$ws.Range("H132").Value = 2989023.5
$ws.Range("I132").Value = 5703009.5
$ws.Range("J132").Value = 3639.3
$ws.Range("K132").Value = 17109028.5
$ws.Range("L132").Value = 10917.9
$ws.Range("M132").Value = -17106498.5
$ws.Range("N132").Value = -15977.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3300.3333
$ws.Range("I132").Value = 3164.9285
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 9494.7855
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -6964.7855
$ws.Range("N132").Value = -15773.4284

$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
